# Helper: replace the full text of paragraph #$paraIndex (1-based) in $textRange
# with $newText while keeping a single run and not disturbing sibling paragraphs'
# formatting. Using Characters(Start,Length) instead of Paragraphs(...).Text avoids
# a spurious run-split on the trailing characters.
function Set-ParaText($textRange, $paraIndex, $newText) {
    $para = $textRange.Paragraphs($paraIndex, 1)
    $full = $textRange.Characters($para.Start, $para.Length)
    $full.Text = $newText
}

$ppres = $ppt.ActivePresentation

# --- Slide 5 ("STGraph - Implementation") -------------------------------
# Content Placeholder 2: tidy up the "Time-Series data layout" bullet list -
# drop stray words / normalize the spacing around the trailing punctuation.
$s5 = $ppres.Slides.Item(5)
$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange

Set-ParaText $tr5 9  "Implemented in AsterixDB ;"
Set-ParaText $tr5 10 "LSM-Tree based;"
Set-ParaText $tr5 11 "Native spatial capabilities;"
Set-ParaText $tr5 12 "Primary index on time dimension;"
Set-ParaText $tr5 13 "Secondary index on spatial dimension."

# --- Slide 6 ("STGraph - Operations") ------------------------------------
$s6 = $ppres.Slides.Item(6)
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1: drop the trailing " -"
Set-ParaText $tr6 1 "Search algorithm: temporal DFS, temporal validity through constraint tightening:"

# Paragraph 2: the formula line also had explicit indent/marL/no-bullet
# overrides on its <a:pPr> that aren't reachable through ParagraphFormat.
# Delete it and re-insert a fresh paragraph (which gets a clean <a:pPr/>)
# right after paragraph 1, then restore its indent level.
$para1 = $tr6.Paragraphs(1, 1)
$para2 = $tr6.Paragraphs(2, 1)
$para2.Delete()
$para1.InsertAfter("`rPath(ni, …, nk) è valido ⇔ ⋂j=i..k-1 I(nj,nj+1) ≠ ∅")
$newPara2 = $tr6.Paragraphs(2, 1)
$newPara2.IndentLevel = 2
